$p = $ppt.ActivePresentation
$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# 1) Footer text on slide 1 and slide 2 (shape #4): bump version/date.
#    " •  labelled version  2.6.0 •  Updated: 2020-09"
#    -> " •  labelled version  2.7.0 •  Updated: 2020-12"
# ---------------------------------------------------------------------------
for ($sIdx = 1; $sIdx -le 2; $sIdx++) {
  $s = $p.Slides.Item($sIdx)
  $sp = $s.Shapes.Item(4)
  $tr = $sp.TextFrame.TextRange
  $t = $tr.Text
  $idx = $t.IndexOf("labelled version  2.6.0")
  if ($idx -ge 0) {
    $runStart = $idx - 4   # back up to include the leading " •  "
    $len = $t.Length - $runStart
    $sub = $tr.Characters($runStart + 1, $len)
    $newText = " " + $bullet + "  labelled version  2.7.0 " + $bullet + "  Updated: 2020-12"
    $sub.Text = $newText
  }
}

# ---------------------------------------------------------------------------
# 2) Slide 2, "Basics" shape (#5): narrower box + retitled question.
#    "When using labelled data?" -> "What is labelled data?"
#    cx 3529812 -> 2955937 EMU
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$basics = $s2.Shapes.Item(5)

$basics.Width = 2955937 / 12700.0

$tr = $basics.TextFrame.TextRange

# "using" (run with err="1") becomes "What" -- keeps that run's formatting.
$rUsing = $tr.Characters(6, 5)
$rUsing.Text = "What"
$tr = $basics.TextFrame.TextRange

# Drop the leading "When ".
$rWhen = $tr.Characters(1, 5)
$rWhen.Text = ""
$tr = $basics.TextFrame.TextRange

# Turn the single space that followed "using" into " is ".
$rSpace = $tr.Characters(5, 1)
$rSpace.Text = " is "
$tr = $basics.TextFrame.TextRange

# Split " is " into " " and "is " (re-assigning with identical text forces a
# run boundary without introducing any new formatting attributes).
$isPart = $tr.Characters(6, 3)
$isPart.Text = "is "

# The shape auto-fits its text box; the diff only touches the width, so pin
# the height back to its original value after the text edits above.
$basics.Height = 340029 / 12700.0

# ---------------------------------------------------------------------------
# 3) Slide 2, shape #18 ("ZoneTexte 67"): fix "identitical" typo and merge
#    the three runs back into the single run the diff expects.
#    "...will not be identitical (original coding will be lost)."
#    -> "...will not be identical (original coding will be lost)."
# ---------------------------------------------------------------------------
$zt67 = $s2.Shapes.Item(18)
$tr = $zt67.TextFrame.TextRange
$t = $tr.Text
$phraseOld = "be identitical (original coding will be lost)."
$phraseNew = "be identical (original coding will be lost)."
$pIdx = $t.IndexOf($phraseOld)
if ($pIdx -ge 0) {
  $sub = $tr.Characters($pIdx + 1, $phraseOld.Length)
  $sub.Text = $phraseNew
}
